$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.687.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.454.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.446.17'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.593'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.67%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.191'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.07'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.562'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000267'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.030.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.472.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.829.77'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '583.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.120'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.850'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '95.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.59'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.79'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.81'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.54'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '572.04'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -20.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.57'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0475'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0956'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.10'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.246.63'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0692'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.295'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '30.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.00%  '
$ws.Range("E47").Value = '  -4.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.126'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  +0.01%  '
